$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $s = $cell.Range.Start
    $len = $oldText.Length
    $rng = $d.Range($s, $s + $len)
    $rng.Text = $newText
}

# Row 1 (cells 1..5)
Set-CellText 1 1 "796÷7=" "125÷4="
Set-CellText 1 2 "821÷6=" "367÷8="
Set-CellText 1 3 "653÷5=" "226÷6="
Set-CellText 1 4 "284÷9=" "397÷3="
Set-CellText 1 5 "543÷2=" "364÷7="

# Row 5 (cells 1..5)
Set-CellText 5 1 "311÷2=" "827÷6="
Set-CellText 5 2 "348÷8=" "695÷4="
Set-CellText 5 3 "709÷2=" "351÷7="
Set-CellText 5 4 "378÷3=" "310÷6="
Set-CellText 5 5 "682÷7=" "660÷2="

# Row 9 (cells 1..5)
Set-CellText 9 1 "466÷8=" "950÷2="
Set-CellText 9 2 "184÷6=" "305÷2="
Set-CellText 9 3 "667÷4=" "792÷7="
Set-CellText 9 4 "291÷9=" "912÷5="
Set-CellText 9 5 "968÷4=" "382÷9="

# Row 13 (cells 1..5)
Set-CellText 13 1 "189÷5=" "524÷5="
Set-CellText 13 2 "165÷2=" "588÷7="
Set-CellText 13 3 "989÷6=" "672÷8="
Set-CellText 13 4 "772÷2=" "993÷4="
Set-CellText 13 5 "914÷6=" "231÷3="

# Row 17 (cells 1..5)
Set-CellText 17 1 "608÷5=" "914÷6="
Set-CellText 17 2 "827÷6=" "759÷3="
Set-CellText 17 3 "688÷8=" "742÷9="
Set-CellText 17 4 "293÷7=" "113÷8="
Set-CellText 17 5 "913÷2=" "539÷6="
